$d = $word.ActiveDocument

# The document currently ends with two empty paragraphs (in this order):
#   ... (O(n)/O(m) paragraph)
#   <empty paragraph>      <- w14:paraId="24AA21A7"
#   <empty paragraph>      <- w14:paraId="15762B0E" (last paragraph, holds sectPr)
#
# We need to insert two new paragraphs ("Oppgave 2", and the explanatory
# paragraph made of three runs) right after the first empty paragraph and
# right before the final empty paragraph.

# Find the empty paragraph that immediately precedes the very last paragraph
# of the document (robust against exact paragraph index assumptions).
$total = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($total - 1)

# Insert a new paragraph after the anchor and fill it with "Oppgave 2".
$anchor.Range.InsertParagraphAfter() | Out-Null
$oppgaveHeading = $d.Paragraphs.Item($total)
$oppgaveHeading.Range.Text = "Oppgave 2"

# Insert another new paragraph after that one; this will hold the three runs.
$oppgaveHeading = $d.Paragraphs.Item($total)
$oppgaveHeading.Range.InsertParagraphAfter() | Out-Null
$bodyPara = $d.Paragraphs.Item($total + 1)

# Build the three runs (separate <w:r> elements, same run formatting) via a
# WordprocessingML package fragment so InsertXML preserves them as distinct
# runs instead of Word's usual "merge adjacent runs with identical rPr"
# behaviour.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>Se javafilene Algorithm1.java og Main1.java for å se algoritme og fungerende test på algoritmefilen i oppgave 1</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> som standard</w:t></w:r>
<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$bodyPara.Range.InsertXML($xml)
